# 242: myr changes (#277)
# Rename "Section 17(4)" to "Vehicle Statistics" and populate its header row
# with the new Vehicle Statistics columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Section 17(4)")
$ws.Name = "Vehicle Statistics"

# Fill in the header row. Columns G (Range) then F (ZEV Type) on purpose,
# so new shared strings are registered in the same order as the source file.
$ws.Cells.Item(1, 1).Value = "Vehicle Class"
$ws.Cells.Item(1, 2).Value = "ZEV Class"
$ws.Cells.Item(1, 3).Value = "Make"
$ws.Cells.Item(1, 4).Value = "Model Name"
$ws.Cells.Item(1, 5).Value = "Model Year"
$ws.Cells.Item(1, 7).Value = "Range"
$ws.Cells.Item(1, 6).Value = "ZEV Type"
$ws.Cells.Item(1, 8).Value = "Submitted Count"
$ws.Cells.Item(1, 9).Value = "Issued Count"

# Bold header row.
$ws.Range("A1:I1").Font.Bold = $true

# Column widths to match the template (values chosen so the saved width,
# after the app's internal rounding, lands on the desired character width).
$ws.Columns.Item(1).ColumnWidth = 15.327
$ws.Columns.Item(2).ColumnWidth = 15.327
$ws.Columns.Item(3).ColumnWidth = 16.165
$ws.Columns.Item(4).ColumnWidth = 17.665
$ws.Columns.Item(5).ColumnWidth = 16.165
$ws.Columns.Item(6).ColumnWidth = 15.165
$ws.Columns.Item(7).ColumnWidth = 19.327
$ws.Columns.Item(8).ColumnWidth = 22.003
$ws.Columns.Item(9).ColumnWidth = 24.003
